$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1 and Q1, copying style from O1 (bold/centered/bordered header style)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update data rows 2-25, columns B through Q (B-O updated values, P/Q new columns)
$vals = 22.71329705539428, 13.79674406896835, 4.541862688152141, 0, 28.49218823937027, 37.0974361257642, 2.35443534782106, 2.939195818955546, 11.95466508217721, 18.35866533933152, 6.202188522811696, 0, 8.37891015590656, 0, 12.67131519049781, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value = $vals[$i] }
$vals = 21.28342497226662, 12.9319893704423, 4.448875291369707, 0, 27.70455149816788, 35.75130169098875, 2.607847720081879, 3.14259692048728, 11.78882243363946, 18.18682259358748, 6.129477573013311, 0, 8.081547706127665, 0, 12.81119916402748, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value = $vals[$i] }
$vals = 20.35439949463182, 12.37434615730549, 4.391198551600767, 0, 27.21834014672674, 34.91233900349651, 2.76888598073543, 3.272532838853643, 11.68950454661704, 18.08491350918565, 6.083229370500147, 0, 7.894080568065824, 0, 12.89877369697971, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value = $vals[$i] }
$vals = 19.95665678583577, 12.14634609618722, 4.368323516267871, 0, 27.00941741734365, 34.54922843096077, 2.836317001273133, 3.329585919846874, 11.64620556899693, 18.03692075052211, 6.063496566681042, 0, 7.817583955879058, 0, 12.93496371315015, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(5, 2 + $i).Value = $vals[$i] }
$vals = 19.88202002407857, 12.11536114800054, 4.36545677076551, 0, 26.96203571342623, 34.46610298443193, 2.848103736950381, 3.34275514796979, 11.63480909650855, 18.02000679306887, 6.059614416803947, 0, 7.80610067051536, 0, 12.94109461152202, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(6, 2 + $i).Value = $vals[$i] }
$vals = 20.32819810001581, 12.39072811107624, 4.393436102055913, 0, 27.18123706551757, 34.84625017723549, 2.771149108376632, 3.282912402086323, 11.67737792091491, 18.05976702412658, 6.081398462435875, 0, 7.896527092075329, 0, 12.8995952075454, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(7, 2 + $i).Value = $vals[$i] }
$vals = 22.20529392602825, 13.52792671052515, 4.513204002159718, 0, 28.17742908048062, 36.55862162229809, 2.442410386748259, 3.02076904488565, 11.88197992880064, 18.26667101377876, 6.175454804746276, 0, 8.281856529548346, 0, 12.71998223889629, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(8, 2 + $i).Value = $vals[$i] }
$vals = 25.51229491499834, 15.51567243151394, 4.738756268144099, 0, 30.15980663397361, 39.90391538179141, 1.837513954901211, 2.53015004465761, 12.32265245818173, 18.74458138188635, 6.350201876645374, 0, 8.989930409098502, 0, 12.37960697444252, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(9, 2 + $i).Value = $vals[$i] }
$vals = 27.61475648473465, 16.87031801279366, 4.877743678442563, 0, 31.39028537121778, 41.99199931113112, 1.692965818429845, 2.754089467060497, 12.60019175117739, 19.00504981831915, 6.452249219989708, 0, 9.384782556584337, 0, 12.1443339820346, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(10, 2 + $i).Value = $vals[$i] }
$vals = 27.89883085999731, 17.60181148009219, 4.739517287819718, 0, 30.32538256243066, 40.65598162935856, 2.699693514146646, 2.845063263517929, 12.26212934797804, 18.2493439437556, 6.389511170488182, 0, 8.728029083210213, 0, 12.12086690143895, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(11, 2 + $i).Value = $vals[$i] }
$vals = 27.72868440280297, 17.92125589494871, 4.600348819976562, 0, 29.19261025464586, 39.14110260922784, 4.030385394575696, 2.853216823290063, 11.92748331605122, 17.58687920256307, 6.391981024415168, 0, 8.096735968766801, 0, 12.17149689628971, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(12, 2 + $i).Value = $vals[$i] }
$vals = 27.168949548254, 17.97568640850428, 4.454842560399648, 0, 27.87070948981967, 37.27604668482658, 5.486700923060134, 2.802607366472842, 11.55471808159942, 16.91726197602613, 6.439681563097209, 0, 7.446019625862993, 0, 12.27817898843101, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(13, 2 + $i).Value = $vals[$i] }
$vals = 26.58779568129002, 17.89962803344607, 4.352546628686481, 0, 26.84325170443821, 35.77858746758933, 6.515180635173873, 2.741987719254877, 11.27300767752241, 16.44180741570733, 6.501375366087042, 0, 6.986238263451909, 0, 12.37927490184455, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(14, 2 + $i).Value = $vals[$i] }
$vals = 26.35463445916021, 17.83015710097409, 4.325428368708051, 0, 26.54296406315325, 35.31913705876691, 6.75456016554797, 2.715817305743961, 11.19335924812933, 16.31796340130296, 6.515923169510016, 0, 6.86965879466922, 0, 12.41360215111958, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(15, 2 + $i).Value = $vals[$i] }
$vals = 25.55285288356949, 17.28119177740631, 4.293278213125236, 0, 26.21663197870872, 34.67422264019957, 6.559230145363524, 2.590741103578607, 11.1352917324608, 16.32289601078034, 6.45745230246605, 0, 6.813450440657588, 0, 12.47931366877192, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(16, 2 + $i).Value = $vals[$i] }
$vals = 25.24965586886876, 16.89363555386313, 4.322578793618468, 0, 26.52581402902372, 35.00246963507795, 5.834291635522081, 2.527093570839611, 11.24262620769742, 16.5788748586589, 6.373630948807342, 0, 7.015850816605259, 0, 12.47757795001949, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(17, 2 + $i).Value = $vals[$i] }
$vals = 25.36637423640663, 16.59450480889744, 4.413928497877373, 0, 27.43706744491035, 36.24318524195104, 4.594521798667697, 2.506630236077981, 11.51083124886788, 17.09945329563192, 6.293755765531429, 0, 7.47912013834776, 0, 12.42461424384476, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(18, 2 + $i).Value = $vals[$i] }
$vals = 25.7767034528906, 16.42309776745891, 4.561383657062164, 0, 28.69957185842909, 38.00547439834549, 3.12484531131104, 2.535180867945264, 11.8676698835721, 17.76083726814029, 6.276287299994481, 0, 8.140917072002352, 0, 12.35286151956159, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(19, 2 + $i).Value = $vals[$i] }
$vals = 27.02481413663207, 16.57101818102579, 4.847463239639852, 0, 30.96670315655725, 41.2752985659724, 1.583272349835602, 2.680774085510864, 12.49084089567489, 18.85880486912986, 6.420761843505065, 0, 9.287608783783217, 0, 12.2095500095996, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(20, 2 + $i).Value = $vals[$i] }
$vals = 28.68103036661392, 17.53626372733509, 4.989445992316793, 0, 32.20447514948984, 43.28021680033161, 1.898818879300171, 2.93299190120335, 12.79519767124238, 19.23234238489399, 6.527033589449749, 0, 9.731054557280643, 0, 12.02066132762369, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(21, 2 + $i).Value = $vals[$i] }
$vals = 29.70166741636563, 18.13420976520153, 5.066089306450496, 0, 32.94899589255414, 44.50221665042557, 2.098058324752444, 3.088351462248175, 12.98272260379547, 19.45939923668418, 6.587551143924052, 0, 9.960737856238882, 0, 11.8999835657265, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(22, 2 + $i).Value = $vals[$i] }
$vals = 29.17918317738012, 17.80041353000824, 5.022593823921846, 0, 32.58659887445389, 43.91048258288799, 1.993347133697995, 3.001750987122444, 12.89481625000868, 19.36473552534792, 6.556922215296626, 0, 9.835288475679709, 0, 11.96267362155374, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(23, 2 + $i).Value = $vals[$i] }
$vals = 27.08233833639354, 16.51627841720043, 4.860604525989193, 0, 31.14620005091086, 41.53955652933907, 1.587366063904803, 2.67489700277749, 12.54654779220969, 18.97218053359023, 6.436355499916588, 0, 9.353445836008861, 0, 12.20518924349686, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(24, 2 + $i).Value = $vals[$i] }
$vals = 24.62972350824339, 15.02987823033017, 4.682619999431839, 0, 29.57093397581762, 38.91416655708912, 1.999024287079091, 2.677427181548386, 12.18155676752859, 18.56932617710409, 6.301610678591703, 0, 8.808924804853833, 0, 12.47200120468252, 0
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(25, 2 + $i).Value = $vals[$i] }
